# "add update to end August" - append new weight-tracker readings (rows 48-70)
# to the raw_data sheet, extend the TOD shared formula down to the new rows,
# and move the sheet selection to the newly-added area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")
$ws.Activate()

# --- New readings (date serial, time serial, weight) for late August 2020 ---
$rows = @(
    @(48, 44074.318055555559, 0.31527777777777777, 74),
    @(49, 44073.915972222225, 0.9159722222222223,  75),
    @(50, 44073.457638888889, 0.45763888888888887, 73.2),
    @(51, 44073.366666666669, 0.3666666666666667,  74.400000000000006),
    @(52, 44072.888888888891, 0.88888888888888884, 75),
    @(53, 44072.381944444445, 0.38194444444444442, 74.599999999999994),
    @(54, 44071.332638888889, 0.33263888888888887, 73.5),
    @(55, 44071.311111111114, 0.31111111111111112, 73.5),
    @(56, 44070.362500000003, 0.36249999999999999, 74.3),
    @(57, 44069.315972222219, 0.31597222222222221, 73.7),
    @(58, 44068.905555555553, 0.90555555555555556, 74.900000000000006),
    @(59, 44068.343055555553, 0.3430555555555555,  74.2),
    @(60, 44068.323611111111, 0.32361111111111113, 74.2),
    @(61, 44068.320138888892, 0.32013888888888892, 74.2),
    @(62, 44067.927777777775, 0.9277777777777777,  73.7),
    @(63, 44067.336111111108, 0.33611111111111108, 74.2),
    @(64, 44067.3125,         0.3125,               74.2),
    @(65, 44066.909722222219, 0.90972222222222221, 75.3),
    @(66, 44066.390277777777, 0.39027777777777778, 73.5),
    @(67, 44064.904861111114, 0.90486111111111101, 74.8),
    @(68, 44064.904166666667, 0.90416666666666667, 74.8),
    @(69, 44064.359722222223, 0.35972222222222222, 73.8),
    @(70, 44074.370833333334, 0.37083333333333335, 73.400000000000006)
)

$lastRow = 70

# Match the date/time number formatting already used lower down the table
# (A19:A47 / B19:B47) by copying their format onto the newly-added rows
# before writing the values.
$ws.Range("A47:B47").Copy() | Out-Null
$ws.Range("A48:B$lastRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# Extend the "AM"/"PM" time-of-day formula down through the new rows. Doing
# this as one bulk range-fill (starting at the existing last formula row, 41)
# makes Excel regroup D41:D70 into its own shared-formula block, matching how
# a fill-down from D41 to D70 behaves.
$ws.Range("D41:D$lastRow").Formula = '=IF(B41<TIME(12,0,0), "AM", "PM")'

# Move the selection/scroll to show the newly-added rows, like the saved view
# in the edited workbook.
$ws.Range("A61").Select()
try {
    $excel.ActiveWindow.ScrollRow = 56
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

$wb.Save()
